$d = $word.ActiveDocument

# Locate the empty paragraph that immediately follows the "Research
# Question" heading - that's the placeholder the draft text belongs in.
$target = $null
for ($i = 2; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $prev = $d.Paragraphs.Item($i - 1)
    if (($p.Range.Text.TrimEnd([char]13, [char]7) -eq "") -and `
        ($prev.Range.Text.TrimEnd([char]13, [char]7) -eq "Research Question")) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the empty paragraph following the 'Research Question' heading."
}

$text1 = "In the modern day games studio artists and designers are often found using keyboard and mouse input to create scenes, art assets and such; for games. However, creative people have a tendency to work better with their hands. The keyboard and mouse input may limit their ability to do this. "
$text2 = "I aim to create a simple tool, where the input is based upon the user in there 3D space as well as using other inputs such as the users voice. Creating an interface more in tune with its user" + [char]0x2019 + "s tendencies. Exploring improvements in productivity and quality of work."

# Fill the empty paragraph with the first chunk of text, split a new
# paragraph after it, then fill that new (still empty) paragraph with the
# second chunk of text.
$r = $target.Range
$r.InsertAfter($text1)
$r.InsertParagraphAfter()
$target.Next().Range.InsertAfter($text2)
